$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# Change 1: split "The general public should have access to naloxone/NARCAN."
# into three runs, bracketing "general public" with proofErr gramStart/gramEnd
# (paragraph 5 - unchanged pPr/rPr, only the run content changes).
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$rPr5 = "<w:rPr><w:rFonts w:ascii='Arial' w:eastAsia='Times New Roman' w:hAnsi='Arial' w:cs='Arial'/><w:color w:val='000000'/></w:rPr>"
$xml5 = "<w:p $wns>" +
        "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='65'/></w:numPr>$rPr5</w:pPr>" +
        "<w:r>$rPr5<w:t xml:space='preserve'>The </w:t></w:r>" +
        "<w:proofErr w:type='gramStart'/>" +
        "<w:r>$rPr5<w:t>general public</w:t></w:r>" +
        "<w:proofErr w:type='gramEnd'/>" +
        "<w:r>$rPr5<w:t xml:space='preserve'> should have access to naloxone/NARCAN.</w:t></w:r>" +
        "</w:p>"
$p5.Range.InsertXML($xml5)

# ---------------------------------------------------------------------------
# Change 2: the "Demographic Items" answer list (numId 66) loses its
# ListParagraph style / numbering and gets explicit left indentation instead
# (360 twips for top-level questions, 1080 twips for the nested answers).
# ---------------------------------------------------------------------------
function Set-PlainIndent($idx, $left) {
    $p = $d.Paragraphs.Item($idx)
    $txt = $p.Range.Text
    $txt = $txt.TrimEnd([char]13, [char]7)
    $escTxt = $txt.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = "<w:p $wns><w:pPr><w:ind w:left='$left'/></w:pPr><w:r><w:t>$escTxt</w:t></w:r></w:p>"
    $p.Range.InsertXML($xml)
}

$simplePairs = @(
    @(48,360),
    @(49,1080),
    @(50,360),
    @(51,1080),
    @(52,1080),
    @(53,1080),
    @(54,1080),
    @(55,1080),
    @(56,360),
    @(57,1080),
    @(58,1080),
    @(59,1080),
    @(60,360),
    @(61,1080),
    @(62,1080),
    @(63,1080),
    @(64,1080),
    @(65,1080),
    @(66,1080),
    @(68,1080),
    @(70,1080),
    @(71,360),
    @(72,1080),
    @(73,1080),
    @(75,1080),
    @(77,1080),
    @(78,1080),
    @(79,1080),
    @(80,1080),
    @(81,360),
    @(82,1080),
    @(83,1080),
    @(84,1080),
    @(85,1080)
)
foreach ($pair in $simplePairs) {
    Set-PlainIndent $pair[0] $pair[1]
}

# Paragraph 67 ("Prefer to self-describe", the race/ethnicity option) also
# picks up the lastRenderedPageBreak marker that used to sit on "Master".
$p67 = $d.Paragraphs.Item(67)
$xml67 = "<w:p $wns><w:pPr><w:ind w:left='1080'/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Prefer to self-describe</w:t></w:r></w:p>"
$p67.Range.InsertXML($xml67)

# Paragraph 69 ("Research requires us to categorize ... below") keeps its two runs.
$p69 = $d.Paragraphs.Item(69)
$xml69 = "<w:p $wns><w:pPr><w:ind w:left='360'/></w:pPr>" +
         "<w:r><w:t>Research requires us to categorize people in racial and ethnic terms</w:t></w:r>" +
         "<w:r><w:t>. Please enter how you would prefer to be described below</w:t></w:r>" +
         "</w:p>"
$p69.Range.InsertXML($xml69)

# Paragraph 74 ("Associate's Degree") keeps its three runs and gains
# proofErr gramStart/gramEnd around the whole phrase.
$p74 = $d.Paragraphs.Item(74)
$xml74 = "<w:p $wns><w:pPr><w:ind w:left='1080'/></w:pPr>" +
         "<w:proofErr w:type='gramStart'/>" +
         "<w:r><w:t>Associate</w:t></w:r>" +
         "<w:r><w:t>’</w:t></w:r>" +
         "<w:r><w:t>s Degree</w:t></w:r>" +
         "<w:proofErr w:type='gramEnd'/>" +
         "</w:p>"
$p74.Range.InsertXML($xml74)

# Paragraph 76 ("Master's Degree") keeps its three runs but loses the
# lastRenderedPageBreak (which moved up to paragraph 67).
$p76 = $d.Paragraphs.Item(76)
$xml76 = "<w:p $wns><w:pPr><w:ind w:left='1080'/></w:pPr>" +
         "<w:r><w:t>Master</w:t></w:r>" +
         "<w:r><w:t>’</w:t></w:r>" +
         "<w:r><w:t>s Degree</w:t></w:r>" +
         "</w:p>"
$p76.Range.InsertXML($xml76)
